$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 1.555284338019139
$ws.Range("E2").Value = 2.194215216927361
$ws.Range("F2").Value = 2.827421753515038
$ws.Range("G2").Value = 3.418112566141071
$ws.Range("H2").Value = 3.947585122819874
$ws.Range("I2").Value = 4.408322322946979
$ws.Range("J2").Value = 4.799168975495887
$ws.Range("K2").Value = 5.122008649631288
$ws.Range("L2").Value = 5.37993129047682
$ws.Range("M2").Value = 5.568630125833018
$ws.Range("N2").Value = 5.692832219368835
$ws.Range("O2").Value = 5.755203733307075
$ws.Range("P2").Value = 5.756275362543696
$ws.Range("Q2").Value = 5.708484311718471
$ws.Range("R2").Value = 5.633708567697529
$ws.Range("S2").Value = 5.546251166966227
$ws.Range("T2").Value = 5.455261312144374
$ws.Range("U2").Value = 5.36639498598111
$ws.Range("V2").Value = 5.282956628988396
$ws.Range("W2").Value = 5.206685113273592
$ws.Range("X2").Value = 5.138294725021634
$ws.Range("Y2").Value = 5.077846634431948
$ws.Range("Z2").Value = 5.025002633391774
$ws.Range("AA2").Value = 4.97919690517586
$ws.Range("AB2").Value = 4.939750693470648
$ws.Range("AC2").Value = 4.905947263102051
$ws.Range("AD2").Value = 4.877079371698912
$ws.Range("AE2").Value = 4.852477860032113
$ws.Range("AF2").Value = 4.83518870344775
